# UniformA-HW45: add the "Holden" HexGrid-90degTilt scheme rows and drop the
# stray duplicated 1Pair-A..MaxUnique header columns (X:AQ) that were left
# over from a copy/paste, shrinking the sheet from A1:AQ19 to A1:W23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the duplicated header columns X:AQ. EntireColumn.Delete shifts
#    everything at/after column X left, leaving rows 1-19 / columns A-W
#    untouched and shrinking the sheet to A1:W19.
$ws.Range("X1:AQ1").EntireColumn.Delete()

# 2. Append four new rows (20-23) for the Holden HexGrid-90degTilt scheme,
#    continuing the existing HKL index sequence (A column) and filling the
#    HKL-combo / pair columns (C:W) with 1s, same as every other row.
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20:W20").Value = 1

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21:W21").Value = 1

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22:W22").Value = 1

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23:W23").Value = 1

# 3. Match the bold/centered/bordered "index" style used by A2:A19 on the
#    new A20:A23 index cells.
$ws.Range("A19").Copy()
$ws.Range("A20:A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
